$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "TestCases" (sheet1): remove the "Login with blank password"
# test case row, and flip RunMode to "Yes" for the newly-enabled cases.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TestCases")

# Row 5 is TC_004 "Login with blank password" - delete it entirely,
# shifting TC_005/TC_006 up.
$ws1.Rows.Item(5).Delete()

# After the delete: row5 = old TC_005 (invalid credentions), row6 = old
# TC_006 (valid credentials). Renumber the literal TestCase_Id values and
# flip RunMode to "Yes" for the newly-enabled cases.
$ws1.Range("C3").Value = "Yes"   # TC_002 Login with blank credentials
$ws1.Range("C4").Value = "Yes"   # TC_003 Login with blank username
$ws1.Range("A5").Value = "TC_004"   # was TC_005, invalid credentions
$ws1.Range("C5").Value = "Yes"
$ws1.Range("A6").Value = "TC_005"   # was TC_006, valid credentials
# Row6 stays RunMode = No

$ws1.Activate()
$ws1.Range("B11").Select()

# ---------------------------------------------------------------------
# Sheet "TestSteps" (sheet2): add the missing steps for TC_003
# (blank -> real password entry) and add brand-new steps for TC_004
# (invalid credentials), plus mark the existing blank TestData values
# with the literal "BLANK" placeholder.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("TestSteps")

# TC_003 "Login with blank username" steps (rows 10-13 already exist blank)
$ws2.Range("A10").Value = "TC_003"
$ws2.Range("B10").Value = "TS_001"
$ws2.Range("C10").Value = "Enter blank username"
$ws2.Range("D10").Value = "sendText"
$ws2.Range("E10").Value = "LOGIN_USERNAME_TEXT_FIELD"

$ws2.Range("A11").Value = "TC_003"
$ws2.Range("B11").Value = "TS_002"
$ws2.Range("C11").Value = "Enter password"
$ws2.Range("D11").Value = "sendText"
$ws2.Range("E11").Value = "LOGIN_PASSWORD_TEXT_FIELD"
$ws2.Range("F11").Value = "qwerty"

$ws2.Range("A12").Value = "TC_003"
$ws2.Range("B12").Value = "TS_003"
$ws2.Range("C12").Value = "Click on login button"
$ws2.Range("D12").Value = "click"
$ws2.Range("E12").Value = "LOGIN_LOGIN_BUTTON"

$ws2.Range("A13").Value = "TC_003"
$ws2.Range("B13").Value = "TS_004"
$ws2.Range("C13").Value = "Verify username validation message"
$ws2.Range("D13").Value = "verifyTextMatches"
$ws2.Range("E13").Value = "LOGIN_USERNAME_ERROR_LABEL"
$ws2.Range("F13").Value = "Required"

# TC_004 "Login with invalid credentions" steps - brand new rows 14-17.
# Seed their formatting from the existing (already-styled) row 13 before
# writing values, so the new cells pick up the same cell style (s="1")
# as every other data row, including the blank trailing G/H cells.
$ws2.Range("A13:H13").Copy()
$ws2.Range("A14:H17").PasteSpecial(-4122)

$ws2.Range("A14").Value = "TC_004"
$ws2.Range("B14").Value = "TS_001"
$ws2.Range("C14").Value = "Enter invalid username"
$ws2.Range("D14").Value = "sendText"
$ws2.Range("E14").Value = "LOGIN_USERNAME_TEXT_FIELD"

$ws2.Range("A15").Value = "TC_004"
$ws2.Range("B15").Value = "TS_002"
$ws2.Range("C15").Value = "Enter invalid password"
$ws2.Range("D15").Value = "sendText"
$ws2.Range("E15").Value = "LOGIN_PASSWORD_TEXT_FIELD"

$ws2.Range("F14").Value = "qweq"
$ws2.Range("F15").Value = "qeqe"

$ws2.Range("A16").Value = "TC_004"
$ws2.Range("B16").Value = "TS_003"
$ws2.Range("C16").Value = "Click on login button"
$ws2.Range("D16").Value = "click"
$ws2.Range("E16").Value = "LOGIN_LOGIN_BUTTON"

$ws2.Range("A17").Value = "TC_004"
$ws2.Range("B17").Value = "TS_004"
$ws2.Range("F17").Value = "Invalid credentials"
$ws2.Range("C17").Value = "Verify invalid credentials validation message"
$ws2.Range("D17").Value = "verifyTextMatches"
$ws2.Range("E17").Value = "LOGIN_INVALID_CREDENTIALS_ERROR_LABEL"

# Mark the blank-field TestData cells with the literal "BLANK" placeholder
# (done last, matching the original authoring order).
$ws2.Range("F5").Value = "BLANK"
$ws2.Range("F6").Value = "BLANK"
$ws2.Range("F10").Value = "BLANK"

$ws2.Activate()
$ws2.Range("A21").Select()

$ws1.Activate()
